$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 145, shifting existing rows 145+ down by one.
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 (duplicate of row 144's data, with an
# updated date and updated volume/price figures).
$ws.Range("A145").Value = 5
$ws.Range("B145").Value = "Macroferia Regional de Talca"
$ws.Range("C145").Value = "Maule"
$ws.Range("D145").Value = 44729
$ws.Range("E145").Value = 7
$ws.Range("F145").Value = 100112008
$ws.Range("G145").Value = "Coliflor"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 3000
$ws.Range("K145").Value = 1000
$ws.Range("L145").Value = 1000
$ws.Range("M145").Value = 1000
$ws.Range("N145").Value = "$/unidad"
$ws.Range("O145").Value = "Región del Maule"
$ws.Range("P145").Value = 1000
$ws.Range("Q145").Value = 1
$ws.Range("R145").Value = "Hortaliza"
